# "simple save/load last query": the Priority of the 'column ambiguously
# defined' feature row is re-classified from HIGH to MED, the tracking
# table is re-sorted by Priority / Est Difficulty (re-running the same
# sort the sheet already remembers), and the last-used cell moves to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Content edit: this feature's Priority moves from HIGH to MED ----
$ws.Range("B2").Value = "MED"

# --- 2. Tidy up column D formatting (a few rows were still on the old
#        11pt / non-wrapped style left over from manual edits) ----------
$dRange = $ws.Range("D2:D17")
$dRange.Font.Size = 10
$dRange.WrapText = $true

# --- 3. Re-apply the table's saved sort (Priority, then Est Difficulty)
#        Only the *visible* (non-filtered) rows participate, exactly like
#        Excel's own Sort does when an AutoFilter is hiding rows. -------
$firstDataRow = 2
$lastDataRow = 17

$snapshot = @()
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowRange = $ws.Rows.Item($r)
    $item = [PSCustomObject]@{
        RowIndex = $r
        Hidden   = $rowRange.Hidden
        Height   = $rowRange.RowHeight
        A        = $ws.Cells.Item($r, 1).Value2
        B        = $ws.Cells.Item($r, 2).Value2
        C        = $ws.Cells.Item($r, 3).Value2
        D        = $ws.Cells.Item($r, 4).Value2
        E        = $ws.Cells.Item($r, 5).Value2
        SortKey  = ""
    }
    $item.SortKey = "{0}|{1}" -f $item.B, $item.C
    $snapshot += $item
}

# Rows hidden by the AutoFilter keep their absolute position; only the
# visible rows are collected and stably re-sorted by Priority then
# Est Difficulty (ascending, text order - matches HIGH < LOW < MED).
$visibleRows = @($snapshot | Where-Object { -not $_.Hidden })
$sortedVisible = @($visibleRows | Sort-Object -Property SortKey)

$finalOrder = New-Object System.Collections.ArrayList
$visibleQueue = New-Object System.Collections.ArrayList
foreach ($row in $sortedVisible) { [void]$visibleQueue.Add($row) }

foreach ($row in $snapshot) {
    if ($row.Hidden) {
        [void]$finalOrder.Add($row)
    } else {
        [void]$finalOrder.Add($visibleQueue[0])
        $visibleQueue.RemoveAt(0)
    }
}

# --- 4. Write the resolved order back into the sheet -------------------
for ($i = 0; $i -lt $finalOrder.Count; $i++) {
    $targetRow = $firstDataRow + $i
    $data = $finalOrder[$i]

    if ($data.RowIndex -ne $targetRow) {
        $ws.Cells.Item($targetRow, 1).Value = $data.A
        $ws.Cells.Item($targetRow, 2).Value = $data.B
        $ws.Cells.Item($targetRow, 3).Value = $data.C
        $ws.Cells.Item($targetRow, 4).Value = $data.D
        $ws.Cells.Item($targetRow, 5).Value = $data.E
        if (-not $ws.Rows.Item($targetRow).Hidden) {
            $ws.Rows.Item($targetRow).RowHeight = $data.Height
        }
    }
}

# --- 5. Last active cell used while reviewing the re-sorted table ------
$ws.Range("D5").Select()
